# ParameterFormulario.xlsx - "Bug: faltaban datos al descargar el excel"
# Adds a second sheet (Hoja1), two workbook-level defined names (Estudios,
# TiposVAlor) and a new "Estudios / Valores" sub-table (rows 19-23) below
# the existing Parameter form, plus capitalises a handful of the
# placeholder tokens that feed the existing form fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameter")

# ---------------------------------------------------------------------------
# 1) Fix-up / re-case the existing merge-field placeholders (values only -
#    formatting of these cells is untouched).
# ---------------------------------------------------------------------------
$ws.Range("B3").Value  = "{{Parameter.Clave}}"
$ws.Range("E3").Value  = "{{Parameter.Areas.Departamento.Nombre}}"
$ws.Range("B5").Value  = "{{Parameter.Nombre}}"
$ws.Range("E5").Value  = "{{Parameter.Areas.Nombre}}"
$ws.Range("B7").Value  = "{{Parameter.NombreCorto}}"
$ws.Range("E7").Value  = "{{Parameter.Reactivos.Nombre}}"
$ws.Range("B9").Value  = "{{Parameter.Unidades}}"
$ws.Range("E9").Value  = "{{Parameter.UnidadSi}}"
$ws.Range("B11").Value = "{{Parameter.TipoValor}}"
$ws.Range("E11").Value = "{{Parameter.Fcsi}}"
$ws.Range("B13").Value = "{{Parameter.Format.Nombre}}"
$ws.Range("E13").Value = "{{Parameter.activo}}"
$ws.Range("B15").Value = "{{Parameter.Formato}}"
$ws.Range("B17").Value = "{{Parameter.ValorInicial}}"

# ---------------------------------------------------------------------------
# 2) New "Estudios" / "Valores" block (rows 19-23, columns A-R)
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "Estudios"
$ws.Range("H19").Value = "Valores"

$ws.Range("A21").Value = "{{item.Id}}"
$ws.Range("B21").Value = "{{item.Nombre}}"
$ws.Range("D21").Value = "nombre"
$ws.Range("E21").Value = "valorInicial"
$ws.Range("F21").Value = "valorFinal"
$ws.Range("G21").Value = "valorInicialNumerico"
$ws.Range("H21").Value = "valorFinalNumerico"
$ws.Range("I21").Value = "rangoEdadInicial"
$ws.Range("J21").Value = "rangoEdadFinal"
$ws.Range("K21").Value = "hombreValorInicial"
$ws.Range("L21").Value = "hombreValorFinal"
$ws.Range("M21").Value = "mujerValorInicial"
$ws.Range("N21").Value = "mujerValorFinal"
$ws.Range("O21").Value = "medidaTiempo"
$ws.Range("P21").Value = "opcion"
$ws.Range("Q21").Value = "descripcionTexto"
$ws.Range("R21").Value = "descripcionParrafo"

$ws.Range("D23").Value = "{{item.nombre}}"
$ws.Range("E23").Value = "{{item.valorInicial}}"
$ws.Range("F23").Value = "{{item.valorFinal}}"
$ws.Range("G23").Value = "{{item.valorInicialNumerico}}"
$ws.Range("H23").Value = "{{item.valorFinalNumerico}}"
$ws.Range("I23").Value = "{{item.rangoEdadInicial}}"
$ws.Range("J23").Value = "{{item.rangoEdadFinal}}"
$ws.Range("K23").Value = "{{item.hombreValorInicial}}"
$ws.Range("L23").Value = "{{item.hombreValorFinal}}"
$ws.Range("M23").Value = "{{item.mujerValorInicial}}"
$ws.Range("N23").Value = "{{item.mujerValorFinal}}"
$ws.Range("O23").Value = "{{item.medidaTiempo}}"
$ws.Range("P23").Value = "{{item.opcion}}"
$ws.Range("Q23").Value = "{{item.descripcionTexto}}"
$ws.Range("R23").Value = "{{item.descripcionParrafo}}"

# Header styling (bold + centered) for the new block - build once on a
# single cell, then fan it out with a formats-only paste so the style
# table only gains the combinations actually used (no partial xfs).
$headerSrc = $ws.Range("A19")
$headerSrc.Font.Bold = $true
$headerSrc.HorizontalAlignment = -4108   # xlCenter
$headerSrc.VerticalAlignment = -4108     # xlCenter
$headerSrc.Copy()
$headerDst = $ws.Range("B19,D19:R20,D21:R21")
$headerDst.PasteSpecial(-4122)           # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Merge the grouped header cells
# ---------------------------------------------------------------------------
$ws.Range("A19:B19").Merge()
$ws.Range("D19:E19").Merge()
$ws.Range("F19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()
$ws.Range("N19:O19").Merge()
$ws.Range("P19:Q19").Merge()

$ws.Range("D20:E20").Merge()
$ws.Range("F20:G20").Merge()
$ws.Range("H20:I20").Merge()
$ws.Range("J20:K20").Merge()
$ws.Range("L20:M20").Merge()
$ws.Range("N20:O20").Merge()
$ws.Range("P20:Q20").Merge()

$ws.Range("D21:E21").Merge()
$ws.Range("F21:G21").Merge()
$ws.Range("H21:I21").Merge()
$ws.Range("J21:K21").Merge()
$ws.Range("L21:M21").Merge()
$ws.Range("N21:O21").Merge()
$ws.Range("P21:Q21").Merge()

# ---------------------------------------------------------------------------
# 4) Column sizing to fit the new wide table
# ---------------------------------------------------------------------------
$ws.Range("A1:R23").Columns.AutoFit()

# ---------------------------------------------------------------------------
# 5) Workbook-level defined names for the new table
# ---------------------------------------------------------------------------
$wb.Names.Add("Estudios", $ws.Range("A21:B22"))
$wb.Names.Add("TiposVAlor", $ws.Range("D23:R24"))

# ---------------------------------------------------------------------------
# 6) Add the trailing empty "Hoja1" worksheet
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$newSheet.Name = "Hoja1"
$ws.Activate()

# Keep the selection where Excel would naturally have left it after this
# editing session.
$ws.Range("E16").Select()
